$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13364.84712035547
$ws.Range("C2").Value = 570778.8208882479

$ws.Range("B3").Value = 21256.78294175362
$ws.Range("C3").Value = 907823.4411594858

$ws.Range("B4").Value = 27475.93170723113
$ws.Range("C4").Value = 1173427.556741277
